# Rename the inline logo pictures in the document's headers/footers.
#
# The Pearson logo (in both footers) goes from "image2.png" to "image1.png",
# and the BTEC logo (in both headers) goes from "image1.jpg" to "image2.jpg".
#
# Renaming is done by selecting the picture first and then updating the
# .Name property via the Selection's InlineShapes collection - addressing
# the InlineShape directly (e.g. $story.Range.InlineShapes(1).Name = ...)
# can fail on some stories with a stale-handle error, but routing the
# assignment through Selection.InlineShapes works reliably for both
# headers and footers.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-LogoPicture($story, $newName) {
    $pic = $story.Range.InlineShapes(1)
    $pic.Select()
    $word.Selection.InlineShapes(1).Name = $newName
}

# Footers: Pearson logo image2.png -> image1.png
Rename-LogoPicture $sec.Footers(1) "image1.png"
Rename-LogoPicture $sec.Footers(2) "image1.png"

# Headers: BTEC logo image1.jpg -> image2.jpg
Rename-LogoPicture $sec.Headers(1) "image2.jpg"
Rename-LogoPicture $sec.Headers(2) "image2.jpg"
